$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M2").Value = "TEST"
$ws.Range("M2").Font.Name = "Aptos Narrow"
